$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.218.99"
$ws.Range("E2").Value = "  +1.10%  "

$ws.Range("D3").Value = "1.882.81"
$ws.Range("E3").Value = "  +1.38%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.007"
$ws.Range("E4").Value = "  +0.33%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "314.39"
$ws.Range("E5").Value = "  +1.08%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.006"
$ws.Range("E6").Value = "  +0.36%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5133"
$ws.Range("E7").Value = "  -0.08%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3920"
$ws.Range("E8").Value = "  +3.09%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.08370"
$ws.Range("E9").Value = "  +1.25%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "1.123"
$ws.Range("E10").Value = "  +1.65%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "6.229"
$ws.Range("E11").Value = "  +1.08%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "20.73"
$ws.Range("E12").Value = "  +1.69%  "

$ws.Range("D13").Value = "1.884.90"
$ws.Range("E13").Value = "  +1.69%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "7.283"
$ws.Range("E14").Value = "  +1.49%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "1.007"
$ws.Range("E15").Value = "  +0.38%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.00001104"
$ws.Range("E16").Value = "  +1.12%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "91.29"
$ws.Range("E17").Value = "  +1.21%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.06662"
$ws.Range("E18").Value = "  +1.05%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "17.82"
$ws.Range("E19").Value = "  +0.90%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "1.005"
$ws.Range("E20").Value = "  +0.30%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.050"
$ws.Range("E21").Value = "  +1.08%  "

$ws.Range("D22").Value = "28.248.80"
$ws.Range("E22").Value = "  +1.07%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "11.17"
$ws.Range("E23").Value = "  +1.59%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.267"
$ws.Range("E24").Value = "  +2.48%  "

$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D25").Value = "2.088.01"
$ws.Range("E25").Value = "  +0.65%  "

$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.510"
$ws.Range("E26").Value = "  -2.14%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "159.43"
$ws.Range("E27").Value = "  +1.83%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "20.63"
$ws.Range("E28").Value = "  +1.53%  "

$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "125.33"
$ws.Range("E29").Value = "  +1.02%  "

$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.1065"
$ws.Range("E30").Value = "  +0.35%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.042"
$ws.Range("E31").Value = "  +0.82%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "5.874"
$ws.Range("E32").Value = "  +5.28%  "

$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.601"
$ws.Range("E33").Value = "  +0.27%  "

$ws.Range("B34").Value = "FraxShare"
$ws.Range("C34").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "9.691"
$ws.Range("E34").Value = "  +2.29%  "

$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.02463"
$ws.Range("E35").Value = "  +2.56%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.06573"
$ws.Range("E36").Value = "  +0.98%  "

$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.2190"
$ws.Range("E37").Value = "  +1.04%  "

$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.208"
$ws.Range("E38").Value = "  +0.41%  "

$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.6523"
$ws.Range("E39").Value = "  +1.71%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.237"
$ws.Range("E40").Value = "  +0.81%  "

$ws.Range("B41").Value = "InternetComputer(DFINITY)"
$ws.Range("C41").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "5.009"
$ws.Range("E41").Value = "  +2.93%  "

$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "11.31"
$ws.Range("E42").Value = "  +0.82%  "

$ws.Range("B43").Value = "Decentraland"
$ws.Range("C43").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.6157"
$ws.Range("E43").Value = "  +1.15%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "13.13"
$ws.Range("E44").Value = "  +0.72%  "

$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.291"
$ws.Range("E45").Value = "  +0.89%  "

$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "3.680"
$ws.Range("E46").Value = "  +0.85%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.015"
$ws.Range("E47").Value = "  +2.51%  "

$ws.Range("B48").Value = "EOS"
$ws.Range("C48").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.237"
$ws.Range("E48").Value = "  +3.11%  "

$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "121.11"
$ws.Range("E49").Value = "  +0.60%  "

$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "78.76"
$ws.Range("E50").Value = "  -0.76%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.06908"
$ws.Range("E51").Value = "  +1.32%  "
